$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (Serial / Code) to row 7
$ws.Range("B7").Value = "FKK128"
$ws.Range("A7").Value = "00000444"

# Update the active selection as recorded in the saved workbook view
$ws.Range("F8:G8").Select()
